# Update code tinh luong % format cac bang
$wb = $excel.ActiveWorkbook

# --- Sheet "Đơn sale chính": drop no-longer-needed columns ---------------
$ws1 = $wb.Worksheets.Item("Đơn sale chính")

# Delete the individual (non-contiguous) columns that are being dropped.
# Must go from right-most to left-most so earlier deletes don't shift the
# still-to-be-deleted column letters out from under us.
$dropCols = @("AA","Y","W","V","U","T","S","R","Q","O","N","I","G")
foreach ($col in $dropCols) {
    $ws1.Range($col + "1").EntireColumn.Delete()
}

# Relabel the headers that remain (after the shift, these sit at G1:N1)
$ws1.Range("G1").Value = "Tên dịch vụ"
$ws1.Range("H1").Value = "Đơn giá gốc"
$ws1.Range("I1").Value = "Sale phụ"
$ws1.Range("J1").Value = "Upsale"
$ws1.Range("K1").Value = "Đơn giá"
$ws1.Range("L1").Value = "Đã thanh toán"
$ws1.Range("M1").Value = "Tỉ lệ chiết khấu sale chính"
$ws1.Range("N1").Value = "Chiết khấu sale chính"

# --- Sheet "Lương": refresh computed salary figures -----------------------
$ws2 = $wb.Worksheets.Item("Lương")

$ws2.Range("A1").Value = "Danh mục lương"
$ws2.Range("B2").Value = 17
$ws2.Range("B3").Value = 595000
$ws2.Range("B12").Value = 6982142.857142857
$ws2.Range("B20").Value = 10473214.28571429
$ws2.Range("B28").Value = -1654000
$ws2.Range("B29").Value = 7892142.857142857
$ws2.Range("B30").Value = 10473214.28571429
$ws2.Range("B31").Value = 16711357.14285714
